$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (values)
$ws.Cells.Item(1, 9).Value2 = "I0"
$ws.Cells.Item(1, 10).Value2 = "IF"

# Copy the header style (bold, bordered, centered) from H1 onto I1 and J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows 2-34 and 37: I = 1, J = same as H
# Rows 35-36: I = 5, J = H + 4
for ($r = 2; $r -le 37; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2

    if ($r -eq 35 -or $r -eq 36) {
        $iVal = 5
        $jVal = $hVal + 4
    } else {
        $iVal = 1
        $jVal = $hVal
    }

    $ws.Cells.Item($r, 9).Value2 = $iVal
    $ws.Cells.Item($r, 10).Value2 = $jVal
}
